$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update MACRO_SCORE column (N) values for rows 2-6 to the new computed value
$newValue = 85.83574689470727
$ws.Range("N2:N6").Value = $newValue
